$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(0.9261335273315012,0.9110007856409936,0.9023844389466831,0.8990431865422579,0.8984986461601068,0.902338689156835,0.9207757203848246,0.9622837274314975,0.9960426473543293,1.012109384398968,1.018295416636278,1.016958613653316,1.012616271374895,1.00996972877391,0.9950069161058366,0.9860093700906702,0.9809010215950593,0.9791829001556209,0.9869602614958239,1.013888957877185,1.032082354421732,1.02231789456107,0.9865301626955443,0.9504816949053634)
$row = 2
foreach ($v in $BValues) {
    $ws.Cells.Item($row, 2).Value2 = $v
    $row++
}

$CValues = @(0.06241763582036697,0.0558126015109508,0.05177497395391129,0.05013399547175368,0.04986177400721203,0.05175282552609417,0.06013643793441759,0.07672445504422853,0.08901100416017016,0.09462410927588394,0.09675322369645301,0.09629452089436086,0.09479920091050076,0.09388374049763115,0.08864466535396787,0.0854368653045583,0.08359406533440961,0.08297050638395831,0.08577810822981746,0.09523831525874016,0.1014418741409031,0.09812897877833393,0.08562382793166989,0.07222018343446734)
$row = 2
foreach ($v in $CValues) {
    $ws.Cells.Item($row, 3).Value2 = $v
    $row++
}

$EValues = @(0.05499493038387882,0.05546402680068008,0.05577272686756585,0.0559037410862615,0.0559258114897685,0.05577447262559732,0.05515239496369384,0.05409573524014744,0.05341783754559781,0.05313059065628067,0.05302483868122243,0.05304748011100302,0.05312182990286285,0.0531677643858941,0.05343703361635832,0.05360762216109904,0.05370773022118702,0.05374196739017045,0.05358925687803406,0.0530999096768161,0.05279770104304227,0.05295738977638287,0.05359755348458961,0.05436422690115261)
$row = 2
foreach ($v in $EValues) {
    $ws.Cells.Item($row, 5).Value2 = $v
    $row++
}

$FValues = @(4.091642545872702,4.050608310526798,4.026585420026947,4.017090569175295,4.015531750792206,4.026456176194472,4.077250545288138,4.186172930012958,4.271906664814196,4.312155916910399,4.327577146183529,4.324247913066444,4.313421028052545,4.306812657193518,4.269301420061367,4.246609394666081,4.233675050469429,4.229315880801209,4.24901283661643,4.316596267936887,4.36181352278976,4.337584312035773,4.24792589288441,4.155706247571743)
$row = 2
foreach ($v in $FValues) {
    $ws.Cells.Item($row, 6).Value2 = $v
    $row++
}

$GValues = @(0.002638425221962101,0.002642672102999894,0.002645417710079299,0.002646571381211121,0.002646765053640892,0.002645433127602982,0.002639860970955642,0.002630023893110314,0.00262345387342756,0.002620606199145135,0.002619548027955417,0.002619775028019484,0.002620518739212068,0.002620976908703227,0.002623642805243298,0.002625314301735166,0.002626288985916214,0.002626621281678799,0.002625134993967215,0.002620299746350642,0.002617257214320337,0.002618870346168925,0.002625216016106293,0.002632569143529658)
$row = 2
foreach ($v in $GValues) {
    $ws.Cells.Item($row, 7).Value2 = $v
    $row++
}

$IValues = @(2.913683916904645,2.883366713492236,2.865333966007455,2.858131559869904,2.856944420300252,2.865236240756474,2.903109531166308,2.982016405214992,3.042849183968286,3.071152077914434,3.081960654427249,3.079628786752778,3.072039482671897,3.06740266036384,3.041012223520681,3.024984073552062,3.015824361893479,3.012733208282853,3.026684162340359,3.074266176056213,3.105893670704688,3.08896488795655,3.025915380632426,2.960170145019262)
$row = 2
foreach ($v in $IValues) {
    $ws.Cells.Item($row, 9).Value2 = $v
    $row++
}

$JValues = @(0.2619274559933054,0.260562304083038,0.2597044850994266,0.2593499441191049,0.2592907712906651,0.2596997238265537,0.2614607965131377,0.2647606071718904,0.2670941136764782,0.2681365577126229,0.2685286051381084,0.2684442903136173,0.2681688656433678,0.2679998091901581,0.2670256070120907,0.266423107404492,0.2660747673065345,0.2659565152530341,0.266487430337186,0.2682498375995479,0.2693859329273209,0.2687810044719612,0.2664583560318441,0.2638840790143888)
$row = 2
foreach ($v in $JValues) {
    $ws.Cells.Item($row, 10).Value2 = $v
    $row++
}

$KValues = @(0.8121239644405307,0.7888830276697263,0.7751941231987303,0.7697618950147387,0.768868702115526,0.7751202706892286,0.8039898549236852,0.8652194110113101,0.913033282487703,0.9354028220721204,0.9439627290057899,0.9421152382901425,0.9361052656954598,0.9324355827708359,0.9115838374804923,0.8989504961981822,0.8917423694050797,0.8893118217703204,0.9002893091633268,0.9378681230134021,0.9629471257898672,0.9495144698898343,0.8996838606690574,0.8481595296674698)
$row = 2
foreach ($v in $KValues) {
    $ws.Cells.Item($row, 11).Value2 = $v
    $row++
}

$LValues = @(0.2551736539138005,0.2557997799285445,0.2563321924790714,0.2565864569789014,0.2566309331411318,0.2563354704184988,0.2553588868680876,0.254614044316952,0.2547755183916962,0.2550019170633249,0.2551095572914619,0.2550854020241999,0.2550103341051511,0.2549672033931216,0.2547637930187747,0.2546781189180791,0.2546432422148541,0.2546339095176364,0.2546857492673951,0.2550317894697827,0.2553856043682856,0.2551851138674479,0.2546822547934795,0.2546907516146746)
$row = 2
foreach ($v in $LValues) {
    $ws.Cells.Item($row, 12).Value2 = $v
    $row++
}

$MValues = @(0.2494508083652391,0.2475364809009051,0.2465320131311621,0.2461657418367196,0.2461075255368002,0.2465268990314442,0.2487553127621815,0.2544789935075897,0.2595071111562319,0.2619728147330918,0.2629321130862969,0.2627243738749954,0.2620512242757513,0.2616422314572411,0.2593495559105428,0.2579887184953193,0.2572227951150978,0.2569663529692754,0.2581318442701743,0.2622482506142774,0.2650876900326509,0.2635586024873575,0.2580670858941119,0.252785905809592)
$row = 2
foreach ($v in $MValues) {
    $ws.Cells.Item($row, 13).Value2 = $v
    $row++
}

Write-Host "applied updates"
